# math_model.xlsx — rework Solver constraints + extend sheet with a new
# "bank credit" helper block (D21:D24, F21) and a new B10/B12/B13 solution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Solver definedNames: the model shrinks from 8 constraints to 4.
#    lhs1..lhs4 now point at B15..B18 (the four computed "period" totals),
#    lhs5..lhs7 collapse onto B18 (same as lhs8 already did).
# ---------------------------------------------------------------------
$wb.Names("solver_lhs1").RefersTo = "=Sheet1!`$B`$15"
$wb.Names("solver_lhs2").RefersTo = "=Sheet1!`$B`$16"
$wb.Names("solver_lhs3").RefersTo = "=Sheet1!`$B`$17"
$wb.Names("solver_lhs4").RefersTo = "=Sheet1!`$B`$18"
$wb.Names("solver_lhs5").RefersTo = "=Sheet1!`$B`$18"
$wb.Names("solver_lhs6").RefersTo = "=Sheet1!`$B`$18"
$wb.Names("solver_lhs7").RefersTo = "=Sheet1!`$B`$18"

# Number of active constraints: 8 -> 4
$wb.Names("solver_num").RefersTo = "=4"

# Relation codes for the 4 active constraints all become "<=" (1)
$wb.Names("solver_rel1").RefersTo = "=1"
$wb.Names("solver_rel2").RefersTo = "=1"
$wb.Names("solver_rel3").RefersTo = "=1"
$wb.Names("solver_rel4").RefersTo = "=1"

# Right-hand sides shift down to match the new lhs1..lhs4 mapping
$wb.Names("solver_rhs1").RefersTo = "=22"
$wb.Names("solver_rhs2").RefersTo = "=25"
$wb.Names("solver_rhs3").RefersTo = "=38"
$wb.Names("solver_rhs4").RefersTo = "=30"
$wb.Names("solver_rhs5").RefersTo = "=30"
$wb.Names("solver_rhs6").RefersTo = "=30"
$wb.Names("solver_rhs7").RefersTo = "=30"

# ---------------------------------------------------------------------
# 2. New Solver solution values for x1..x4 (B10, B12, B13 change; B11
#    stays 0). Downstream formulas (F10, B15:B18) recompute automatically.
# ---------------------------------------------------------------------
$ws.Range("B10").Value = 2.4210524519575838
$ws.Range("B12").Value = 0.52631597972048771
$ws.Range("B13").Value = 0

# Row 17's formula coefficient changed from 19 to 10
$ws.Range("B17").Formula = "= 10 * B10 + 9 * B11 + 9 * B12 + 7 * B13"

# ---------------------------------------------------------------------
# 3. New header label in H1 + a small "bank credit" check block below
#    the model (D21:D24 values, F21 boolean constraint check).
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "количество проектов претендующий на получение кредита в банке"

$ws.Range("D21").Value = -0.8
$ws.Range("D22").Value = -0.98
$ws.Range("D23").Value = 1.9
$ws.Range("D24").Value = 0.4
$ws.Range("F21").Formula = "= 8 * D21 + 7 *D22 + 5 * D23 + 9 * D24 <= 22"

# ---------------------------------------------------------------------
# 4. Selection moves to F14 (matches the saved cursor position).
# ---------------------------------------------------------------------
$null = $ws.Range("F14").Select()
